$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A362").Value = 'Obed Wild & Scenic River'
$ws.Range("A363").Value = 'Ocmulgee Mounds National Historical Park'
$ws.Range("A364").Value = 'Oil Region National Heritage Area'
$ws.Range("A365").Value = 'Oklahoma City National Memorial'
$ws.Range("A366").Value = 'Old Spanish National Historic Trail'
$ws.Range("A367").Value = 'Olympic National Park'
$ws.Range("A368").Value = 'Oregon Caves National Monument & Preserve'
$ws.Range("A369").Value = 'Oregon National Historic Trail'
$ws.Range("A370").Value = 'Organ Pipe Cactus National Monument'
$ws.Range("A371").Value = 'Overmountain Victory National Historic Trail'
$ws.Range("A372").Value = 'Oxon Cove  Park & Oxon Hill Farm'
$ws.Range("A373").Value = 'Ozark National Scenic Riverways'
$ws.Range("A375").Value = 'Padre Island National Seashore'
$ws.Range("A376").Value = 'Palo Alto Battlefield National Historical Park'
$ws.Range("A377").Value = 'Parashant National Monument'
$ws.Range("A378").Value = 'Paterson Great Falls National Historical Park'
$ws.Range("A379").Value = 'Pea Ridge National Military Park'
$ws.Range("A380").Value = 'Pearl Harbor National Memorial'
$ws.Range("A381").Value = 'Pecos National Historical Park'
$ws.Range("A382").Value = 'Pennsylvania Avenue'
$ws.Range("A383").Value = 'Perry''s Victory & International Peace Memorial'
$ws.Range("A384").Value = 'Petersburg National Battlefield'
$ws.Range("A385").Value = 'Petrified Forest National Park'
$ws.Range("A386").Value = 'Petroglyph National Monument'
$ws.Range("A387").Value = 'Pictured Rocks National Lakeshore'
$ws.Range("A388").Value = 'Pinnacles National Park'
$ws.Range("A389").Value = 'Pipe Spring National Monument'
$ws.Range("A390").Value = 'Pipestone National Monument'
$ws.Range("A391").Value = 'Piscataway Park'
$ws.Range("A392").Value = 'Point Reyes National Seashore'
$ws.Range("A393").Value = 'Pony Express National Historic Trail'
$ws.Range("A394").Value = 'Port Chicago Naval Magazine National Memorial'
$ws.Range("A395").Value = 'Potomac Heritage National Scenic Trail'
$ws.Range("A396").Value = 'Poverty Point National Monument'
$ws.Range("A397").Value = 'President William Jefferson Clinton Birthplace Home National Historic Site'
$ws.Range("A398").Value = 'President''s Park'
$ws.Range("A399").Value = 'Presidio of San Francisco'
$ws.Range("A400").Value = 'Prince William Forest Park'
$ws.Range("A401").Value = 'Pu`uhonua O Hōnaunau National Historical Park'
$ws.Range("A402").Value = 'Pu`ukoholā Heiau National Historic Site'
$ws.Range("A403").Value = 'Pullman National Monument'
$ws.Range("A405").Value = 'Rainbow Bridge National Monument'
$ws.Range("A406").Value = 'Reconstruction Era National Historical Park'
$ws.Range("A407").Value = 'Redwood National and State Parks'
$ws.Range("A408").Value = 'Richmond National Battlefield Park'
$ws.Range("A409").Value = 'Rio Grande Wild & Scenic River'
$ws.Range("A410").Value = 'River Raisin National Battlefield Park'
$ws.Range("A411").Value = 'Rivers Of Steel National Heritage Area'
$ws.Range("A412").Value = 'Rock Creek Park'
$ws.Range("A413").Value = 'Rocky Mountain National Park'
$ws.Range("A414").Value = 'Roger Williams National Memorial'
$ws.Range("A415").Value = 'Roosevelt Campobello International Park'
$ws.Range("A416").Value = 'Rosie the Riveter WWII Home Front National Historical Park'
$ws.Range("A417").Value = 'Russell Cave National Monument'
$ws.Range("A419").Value = 'Sagamore Hill National Historic Site'
$ws.Range("A420").Value = 'Saguaro National Park'
$ws.Range("A421").Value = 'Saint Croix Island International Historic Site'
$ws.Range("A422").Value = 'Saint Croix National Scenic Riverway'
$ws.Range("A423").Value = 'Saint Paul''s Church National Historic Site'
$ws.Range("A424").Value = 'Saint-Gaudens National Historical Park'
$ws.Range("A425").Value = 'Salem Maritime National Historic Site'
$ws.Range("A426").Value = 'Salinas Pueblo Missions National Monument'
$ws.Range("A427").Value = 'Salt River Bay National Historical Park and Ecological Preserve'
$ws.Range("A428").Value = 'San Antonio Missions National Historical Park'
$ws.Range("A429").Value = 'San Francisco Maritime National Historical Park'
$ws.Range("A430").Value = 'San Juan Island National Historical Park'
$ws.Range("A431").Value = 'San Juan National Historic Site'
$ws.Range("A432").Value = 'Sand Creek Massacre National Historic Site'
$ws.Range("A433").Value = 'Santa Fe National Historic Trail'
$ws.Range("A434").Value = 'Santa Monica Mountains National Recreation Area'
$ws.Range("A435").Value = 'Saratoga National Historical Park'
$ws.Range("A436").Value = 'Saugus Iron Works National Historic Site'
$ws.Range("A437").Value = 'Schuylkill River Valley National Heritage Area'
$ws.Range("A438").Value = 'Scotts Bluff National Monument'
$ws.Range("A439").Value = 'Selma To Montgomery National Historic Trail'
$ws.Range("A440").Value = 'Sequoia & Kings Canyon National Parks'
$ws.Range("A441").Value = 'Shenandoah National Park'
$ws.Range("A442").Value = 'Shenandoah Valley Battlefields National Historic District'
$ws.Range("A443").Value = 'Shiloh National Military Park'
$ws.Range("A444").Value = 'Sitka National Historical Park'
$ws.Range("A445").Value = 'Sleeping Bear Dunes National Lakeshore'
$ws.Range("A446").Value = 'South Carolina National Heritage Corridor'
$ws.Range("A447").Value = 'Springfield Armory National Historic Site'
$ws.Range("A448").Value = 'Star-Spangled Banner National Historic Trail'
$ws.Range("A449").Value = 'Statue Of Liberty National Monument'
$ws.Range("A450").Value = 'Steamtown National Historic Site'
$ws.Range("A451").Value = 'Stones River National Battlefield'
$ws.Range("A452").Value = 'Stonewall National Monument'
$ws.Range("A453").Value = 'Sunset Crater Volcano National Monument'
$ws.Range("A455").Value = 'Tallgrass Prairie National Preserve'
$ws.Range("A456").Value = 'Tennessee Civil War National Heritage Area'
$ws.Range("A457").Value = 'Thaddeus Kosciuszko National Memorial'
$ws.Range("A458").Value = 'The Last Green Valley National Heritage Corridor'
$ws.Range("A459").Value = 'Theodore Roosevelt Birthplace National Historic Site'
$ws.Range("A460").Value = 'Theodore Roosevelt Inaugural National Historic Site'
$ws.Range("A461").Value = 'Theodore Roosevelt Island'
$ws.Range("A462").Value = 'Theodore Roosevelt National Park'
$ws.Range("A463").Value = 'Thomas Cole National Historic Site'
$ws.Range("A464").Value = 'Thomas Edison National Historical Park'
$ws.Range("A465").Value = 'Thomas Jefferson Memorial'
$ws.Range("A466").Value = 'Thomas Stone National Historic Site'
$ws.Range("A467").Value = 'Timpanogos Cave National Monument'
$ws.Range("A468").Value = 'Timucuan Ecological & Historic Preserve'
$ws.Range("A469").Value = 'Tonto National Monument'
$ws.Range("A470").Value = 'Touro Synagogue National Historic Site'
$ws.Range("A471").Value = 'Trail Of Tears National Historic Trail'
$ws.Range("A472").Value = 'Tule Lake National Monument'
$ws.Range("A473").Value = 'Tule Springs Fossil Beds National Monument'
$ws.Range("A474").Value = 'Tumacácori National Historical Park'
$ws.Range("A475").Value = 'Tupelo National Battlefield'
$ws.Range("A476").Value = 'Tuskegee Airmen National Historic Site'
$ws.Range("A477").Value = 'Tuskegee Institute National Historic Site'
$ws.Range("A478").Value = 'Tuzigoot National Monument'
$ws.Range("A480").Value = 'Ulysses S Grant National Historic Site'
$ws.Range("A481").Value = 'Upper Delaware Scenic & Recreational River'
$ws.Range("A482").Value = 'Upper Housatonic Valley National Heritage Area'
$ws.Range("A484").Value = 'Valles Caldera National Preserve'
$ws.Range("A485").Value = 'Valley Forge National Historical Park'
$ws.Range("A486").Value = 'Vanderbilt Mansion National Historic Site'
$ws.Range("A487").Value = 'Vicksburg National Military Park'
$ws.Range("A488").Value = 'Vietnam Veterans Memorial'
$ws.Range("A489").Value = 'Virgin Islands Coral Reef National Monument'
$ws.Range("A490").Value = 'Virgin Islands National Park'
$ws.Range("A491").Value = 'Voyageurs National Park'
$ws.Range("A493").Value = 'Waco Mammoth National Monument'
$ws.Range("A494").Value = 'Walnut Canyon National Monument'
$ws.Range("A495").Value = 'War In The Pacific National Historical Park'
$ws.Range("A496").Value = 'Washington Monument'
$ws.Range("A497").Value = 'Washington-Rochambeau Revolutionary Route National Historic Trail'
$ws.Range("A498").Value = 'Washita Battlefield National Historic Site'
$ws.Range("A499").Value = 'Weir Farm National Historic Site'
$ws.Range("A500").Value = 'Wheeling National Heritage Area'
$ws.Range("A501").Value = 'Whiskeytown National Recreation Area'
$ws.Range("A502").Value = 'White Sands National Monument'
$ws.Range("A503").Value = 'Whitman Mission National Historic Site'
$ws.Range("A504").Value = 'William Howard Taft National Historic Site'
$ws.Range("A505").Value = 'Wilson''s Creek National Battlefield'
$ws.Range("A506").Value = 'Wind Cave National Park'
$ws.Range("A507").Value = 'Wing Luke Museum Affiliated Area'
$ws.Range("A508").Value = 'Wolf Trap National Park for the Performing Arts'
$ws.Range("A509").Value = 'Women''s Rights National Historical Park'
$ws.Range("A510").Value = 'World War II Memorial'
$ws.Range("A511").Value = 'Wrangell - St Elias National Park & Preserve'
$ws.Range("A512").Value = 'Wright Brothers National Memorial'
$ws.Range("A513").Value = 'Wupatki National Monument'
$ws.Range("A515").Value = 'Yellowstone National Park'
$ws.Range("A516").Value = 'Yorktown Battlefield Part of Colonial National Historical Park'
$ws.Range("A517").Value = 'Yosemite National Park'
$ws.Range("A518").Value = 'Yucca House National Monument'
$ws.Range("A519").Value = 'Yukon - Charley Rivers National Preserve'
$ws.Range("A521").Value = 'Zion National Park'

$ws.Range("A520").Select()
